$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated typology mapping text (85:15 CDN/CDL split) ---

$officesText = @"
19.21% CR/LFINF+CDN/HBET:3-5/Offices
3.39% CR/LFINF+CDL/HBET:3-5/Offices
2.125% CR/LFINF+CDN/HBET:6-/Offices
0.375% CR/LFINF+CDL/HBET:6-/Offices
24.82% CR/LWAL+CDN/HBET:3-5/Offices
4.38% CR/LWAL+CDL/HBET:3-5/Offices
2.72% CR/LWAL+CDN/HBET:6-/Offices
0.48% CR/LWAL+CDL/HBET:6-/Offices
1.0% MUR+CB/LWAL+CDN/H:1/Offices
2.3% MUR+CB/LWAL+CDN/H:2/Offices
9.1% MUR+CL/LWAL+CDN/H:1/Offices
21.1% MUR+CL/LWAL+CDN/H:2/Offices
2.7% MUR+ST/LWAL+CDN/H:1/Offices
6.3% MUR+ST/LWAL+CDN/H:2/Offices
0.0% CR/LFINF+CDL/H:1/Offices
0.0% CR/LFINF+CDL/H:2/Offices
0.0% W/LPB+CDL/H:1/Offices
0.0% W/LPB+CDL/H:2/Offices
"@

$tradeText = @"
0.0% CR/LFINF+CDL/HBET:3-5/Trade
0.0% CR/LFINF+CDL/HBET:6-/Trade
0.0% CR/LWAL+CDL/HBET:3-5/Trade
0.0% CR/LWAL+CDL/HBET:6-/Trade
2.3% MUR+CB/LWAL+CDN/H:1/Trade
1.0% MUR+CB/LWAL+CDN/H:2/Trade
21.1% MUR+CL/LWAL+CDN/H:1/Trade
9.1% MUR+CL/LWAL+CDN/H:2/Trade
6.3% MUR+ST/LWAL+CDN/H:1/Trade
2.7% MUR+ST/LWAL+CDN/H:2/Trade
7.725% CR/LFINF+CDL/H:1/Trade
0.855% CR/LFINF+CDL/H:2/Trade
43.775% CR/LFINF+CDN/H:1/Trade
4.845% CR/LFINF+CDN/H:2/Trade
0.3% W/LPB+CDL/H:1/Trade
0.0% W/LPB+CDL/H:2/Trade
"@

$hotelsText = @"
3.36% CR/LFINF+CDL/HBET:3-5/Hotels
0.375% CR/LFINF+CDL/HBET:6-/Hotels
4.365% CR/LWAL+CDL/HBET:3-5/Hotels
0.48% CR/LWAL+CDL/HBET:6-/Hotels
19.04% CR/LFINF+CDN/HBET:3-5/Hotels
2.125% CR/LFINF+CDN/HBET:6-/Hotels
24.735% CR/LWAL+CDN/HBET:3-5/Hotels
2.72% CR/LWAL+CDN/HBET:6-/Hotels
 1.0% MUR+CB/LWAL+CDN/H:1/Hotels
 2.3% MUR+CB/LWAL+CDN/H:2/Hotels
 21.1% MUR+CL/LWAL+CDN/H:1/Hotels
 9.1% MUR+CL/LWAL+CDN/H:2/Hotels
 6.3% MUR+ST/LWAL+CDN/H:1/Hotels
 2.7% MUR+ST/LWAL+CDN/H:2/Hotels
 0.0% CR/LFINF+CDL/H:1/Hotels
 0.0% CR/LFINF+CDL/H:2/Hotels
 0.3% W/LPB+CDL/H:1/Hotels
 0.0% W/LPB+CDL/H:2/Hotels
"@

$ws.Range("B2").Value = $officesText.TrimEnd("`r", "`n")
$ws.Range("C2").Value = $tradeText.TrimEnd("`r", "`n")
$ws.Range("D2").Value = $hotelsText.TrimEnd("`r", "`n")

# --- Formatting: wrap the long mapping text, widen the columns, and
#     grow row 2 to fit the extra lines ---

$ws.Range("B2:D2").WrapText = $true

$ws.Columns.Item(2).ColumnWidth = 31.166666666666668
$ws.Columns.Item(3).ColumnWidth = 32.998697916666664
$ws.Columns.Item(4).ColumnWidth = 31.498697916666668

$ws.Rows.Item(2).RowHeight = 289

# --- Move the active selection, matching the saved workbook state ---
$ws.Range("F2").Select()
